$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1093.7812
$ws.Range("I15").Value = 1093.7812
$ws.Range("K15").Value = 3281.3436
$ws.Range("M15").Value = -3112.3436

$ws.Range("H99").Value = 1063.2858
$ws.Range("I99").Value = 1061
$ws.Range("J99").Value = 1066.3334
$ws.Range("K99").Value = 3183
$ws.Range("L99").Value = 3199.0002
$ws.Range("M99").Value = -1685
$ws.Range("N99").Value = -6195.0002

$ws.Range("H107").Value = 407.77777
$ws.Range("I107").Value = 414.2857
$ws.Range("J107").Value = 385
$ws.Range("K107").Value = 414.2857
$ws.Range("L107").Value = 385
$ws.Range("M107").Value = 1505.7143
$ws.Range("N107").Value = -4225

$ws.Range("H116").Value = 7526.6
$ws.Range("I116").Value = 6067.6665
$ws.Range("K116").Value = 6067.6665
$ws.Range("M116").Value = -2625.6665

$ws.Range("H137").Value = 2532.6667
$ws.Range("I137").Value = 2196.2
$ws.Range("K137").Value = 6588.599999999999
$ws.Range("M137").Value = -4038.599999999999

$ws.Range("H141").Value = 1959.75
$ws.Range("I141").Value = 1959.75
$ws.Range("K141").Value = 5879.25
$ws.Range("M141").Value = -699.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11917.134
$ws.Range("I32").Value = 10673.692
$ws.Range("J32").Value = 19999.5
$ws.Range("K32").Value = 10673.692
$ws.Range("L32").Value = 19999.5
$ws.Range("M32").Value = -10386.692
$ws.Range("N32").Value = -20573.5

$ws.Range("H74").Value = 2612.1667
$ws.Range("I74").Value = 1949.1333
$ws.Range("J74").Value = 3717.2222
$ws.Range("K74").Value = 1949.1333
$ws.Range("L74").Value = 3717.2222
$ws.Range("M74").Value = -1075.1333
$ws.Range("N74").Value = -5465.2222

$ws.Range("H77").Value = 2612.1667
$ws.Range("I77").Value = 1949.1333
$ws.Range("J77").Value = 3717.2222
$ws.Range("K77").Value = 9745.666499999999
$ws.Range("L77").Value = 18586.111
$ws.Range("M77").Value = -5377.666499999999
$ws.Range("N77").Value = -27322.111

$ws.Range("H102").Value = 2119.5
$ws.Range("I102").Value = 2119.5
$ws.Range("K102").Value = 2119.5
$ws.Range("M102").Value = -497.5

$ws.Range("H110").Value = 2716
$ws.Range("I110").Value = 2716
$ws.Range("K110").Value = 2716
$ws.Range("M110").Value = -671

$ws.Range("H122").Value = 2020.6923
$ws.Range("I122").Value = 1531.7142
$ws.Range("K122").Value = 4595.142599999999
$ws.Range("M122").Value = -2145.142599999999

$ws.Range("H127").Value = 59999
$ws.Range("J127").Value = 59999
$ws.Range("L127").Value = 59999
$ws.Range("N127").Value = -69919

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 30994
$ws.Range("J6").Value = 30994
$ws.Range("L6").Value = 30994
$ws.Range("N6").Value = -31220

$ws.Range("H25").Value = 6425
$ws.Range("I25").Value = 7433.3335
$ws.Range("J25").Value = 3400
$ws.Range("K25").Value = 7433.3335
$ws.Range("L25").Value = 3400
$ws.Range("M25").Value = -7198.3335
$ws.Range("N25").Value = -3870

$ws.Range("H82").Value = 22499.6
$ws.Range("I82").Value = 9374.75
$ws.Range("K82").Value = 9374.75
$ws.Range("M82").Value = -8991.75

$ws.Range("H85").Value = 22499.6
$ws.Range("I85").Value = 9374.75
$ws.Range("K85").Value = 9374.75
$ws.Range("M85").Value = -8048.75

$ws.Range("H134").Value = 1741
$ws.Range("I134").Value = 1741
$ws.Range("K134").Value = 5223
$ws.Range("M134").Value = -2688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1711.9445
$ws.Range("I16").Value = 1208.0769
$ws.Range("J16").Value = 3022
$ws.Range("K16").Value = 1208.0769
$ws.Range("L16").Value = 3022
$ws.Range("M16").Value = -921.0769
$ws.Range("N16").Value = -3596

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0

$ws.Range("H107").Value = 909.6
$ws.Range("I107").Value = 530.5
$ws.Range("K107").Value = 530.5
$ws.Range("M107").Value = 1389.5

$ws.Range("H113").Value = 1711.9445
$ws.Range("I113").Value = 1208.0769
$ws.Range("J113").Value = 3022
$ws.Range("K113").Value = 1208.0769
$ws.Range("L113").Value = 3022
$ws.Range("M113").Value = 961.9231
$ws.Range("N113").Value = -7362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38.125
$ws.Range("J12").Value = 16.4
$ws.Range("L12").Value = 49.2
$ws.Range("N12").Value = -395.2

$ws.Range("H68").Value = 4557.6
$ws.Range("I68").Value = 484
$ws.Range("J68").Value = 5576
$ws.Range("K68").Value = 1452
$ws.Range("L68").Value = 16728
$ws.Range("M68").Value = -641
$ws.Range("N68").Value = -18350

$ws.Range("H71").Value = 4557.6
$ws.Range("I71").Value = 484
$ws.Range("J71").Value = 5576
$ws.Range("K71").Value = 4356
$ws.Range("L71").Value = 50184
$ws.Range("M71").Value = -300
$ws.Range("N71").Value = -58296

$ws.Range("H130").Value = 2581.6667
$ws.Range("I130").Value = 2003.3334
$ws.Range("J130").Value = 3160
$ws.Range("K130").Value = 6010.0002
$ws.Range("L130").Value = 9480
$ws.Range("M130").Value = -990.0002000000004
$ws.Range("N130").Value = -19520

$ws.Range("H140").Value = 6544.263
$ws.Range("I140").Value = 2641.8462
$ws.Range("K140").Value = 7925.5386
$ws.Range("M140").Value = -2745.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 8000
$ws.Range("I46").Value = 8000
$ws.Range("K46").Value = 8000
$ws.Range("M46").Value = -7844

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0

$ws.Range("H80").Value = 2887.5293
$ws.Range("I80").Value = 1658.75
$ws.Range("K80").Value = 1658.75
$ws.Range("M80").Value = -660.75

$ws.Range("H83").Value = 2887.5293
$ws.Range("I83").Value = 1658.75
$ws.Range("K83").Value = 8293.75
$ws.Range("M83").Value = -3301.75

$ws.Range("H113").Value = 1176.4445
$ws.Range("I113").Value = 1132.3334
$ws.Range("K113").Value = 1132.3334
$ws.Range("M113").Value = 1037.6666

$ws.Range("H122").Value = 5581.2
$ws.Range("I122").Value = 4400
$ws.Range("J122").Value = 5876.5
$ws.Range("K122").Value = 13200
$ws.Range("L122").Value = 17629.5
$ws.Range("M122").Value = -10750
$ws.Range("N122").Value = -22529.5

$ws.Range("H132").Value = 3058.25
$ws.Range("I132").Value = 3058.25
$ws.Range("K132").Value = 9174.75
$ws.Range("M132").Value = -6644.75

$ws.Range("H136").Value = 28656
$ws.Range("J136").Value = 28656
$ws.Range("L136").Value = 85968
$ws.Range("N136").Value = -91068

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 420.5
$ws.Range("I16").Value = 421.45456
$ws.Range("K16").Value = 421.45456
$ws.Range("M16").Value = -251.45456

$ws.Range("H61").Value = 3131.7778
$ws.Range("I61").Value = 2433.7646
$ws.Range("K61").Value = 2433.7646
$ws.Range("M61").Value = -2231.7646

$ws.Range("H93").Value = 3179.5
$ws.Range("I93").Value = 2907.4
$ws.Range("J93").Value = 3633
$ws.Range("K93").Value = 2907.4
$ws.Range("L93").Value = 3633
$ws.Range("M93").Value = -1659.4
$ws.Range("N93").Value = -6129

$ws.Range("H113").Value = 3131.7778
$ws.Range("I113").Value = 2433.7646
$ws.Range("K113").Value = 2433.7646
$ws.Range("M113").Value = -263.7646

$ws.Range("H122").Value = 7287.6787
$ws.Range("I122").Value = 8129.0625
$ws.Range("K122").Value = 24387.1875
$ws.Range("M122").Value = -21937.1875

$ws.Range("H132").Value = 2576.6667
$ws.Range("I132").Value = 1238.75
$ws.Range("K132").Value = 3716.25
$ws.Range("M132").Value = -1186.25

$ws.Range("H136").Value = 4961.875
$ws.Range("I136").Value = 4978.2
$ws.Range("J136").Value = 4934.6665
$ws.Range("K136").Value = 14934.6
$ws.Range("L136").Value = 14803.9995
$ws.Range("M136").Value = -12384.6
$ws.Range("N136").Value = -19903.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 200000000
$ws.Range("I49").Value = 200000000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 200000000
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("M49").Value = -199999770

$ws.Range("H75").Value = 88559
$ws.Range("I75").Value = 88559
$ws.Range("K75").Value = 88559
$ws.Range("M75").Value = -87623

$ws.Range("H78").Value = 88559
$ws.Range("I78").Value = 88559
$ws.Range("K78").Value = 265677
$ws.Range("M78").Value = -260997
